# M2 - Experiment.xlsx edit script
# Updates the raw timing-gate data in the "medium" and "long" sheets
# (columns A/B, rows 2-42) with re-measured values, and restores the
# expected sheet selection / active-tab state ("short" becomes the
# active tab, with C1:D42 selected; "medium" ends up with A1:B42
# selected but is not the active tab).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. "medium" sheet - replace A2:B42 with the corrected data set
# ---------------------------------------------------------------
$medium = $wb.Worksheets.Item("medium")

$mediumData = @(
    @(2, 1326, 9.3999999999999986),
    @(3, 1334.5, 9.3000000000000007),
    @(4, 1332, 9.6999999999999993),
    @(5, 1335.5, 9.6),
    @(6, 1342, 8.9),
    @(7, 1339.5, 5.9),
    @(8, 1344, 6.5),
    @(9, 1345, 7.1000000000000005),
    @(10, 1346, 11.5),
    @(11, 1348.5, 14.9),
    @(12, 1345, 17.3),
    @(13, 1350, 21.7),
    @(14, 1347, 26.9),
    @(15, 1353, 28.8),
    @(16, 1350, 28.8),
    @(17, 1352.5, 25.3),
    @(18, 1351, 29.8),
    @(19, 1353, 28.7),
    @(20, 1351, 27.900000000000002),
    @(21, 1359.5, 24.3),
    @(22, 1356, 21.2),
    @(23, 1358.5, 19.7),
    @(24, 1363, 15.8),
    @(25, 1368.5, 14.5),
    @(26, 1366, 11.299999999999999),
    @(27, 1371.5, 11.5),
    @(28, 1371, 8.5),
    @(29, 1379.5, 11.399999999999999),
    @(30, 1382, 12.7),
    @(31, 1378.5, 15.6),
    @(32, 1385, 21.5),
    @(33, 1385.5, 26.1),
    @(34, 1390, 25.1),
    @(35, 1390.5, 22.400000000000002),
    @(36, 1397, 17.600000000000001),
    @(37, 1393.5, 14.799999999999999),
    @(38, 1398, 11.2),
    @(39, 1399.5, 10.1),
    @(40, 1401, 10.1),
    @(41, 1406.5, 10.200000000000001),
    @(42, 1409, 19)
)

foreach ($row in $mediumData) {
    $r = $row[0]
    $medium.Cells.Item($r, 1).Value = $row[1]
    $medium.Cells.Item($r, 2).Value = $row[2]
}

# ---------------------------------------------------------------
# 2. "long" sheet - replace A2:B42 with the corrected data set
# ---------------------------------------------------------------
$long = $wb.Worksheets.Item("long")

$longData = @(
    @(2, 1317, 9.3999999999999986),
    @(3, 1315.5, 8.8000000000000007),
    @(4, 1312, 9.6),
    @(5, 1323.5, 9.2999999999999989),
    @(6, 1328, 8.3000000000000007),
    @(7, 1325.5, 6.4),
    @(8, 1332, 6.7),
    @(9, 1325, 6.7),
    @(10, 1335, 12.3),
    @(11, 1330.5, 14.3),
    @(12, 1327, 17.700000000000003),
    @(13, 1335, 21.2),
    @(14, 1332, 27.5),
    @(15, 1333, 28.400000000000002),
    @(16, 1335, 28.7),
    @(17, 1329.5, 26),
    @(18, 1334, 29.7),
    @(19, 1332, 28.7),
    @(20, 1344, 27.5),
    @(21, 1336.5, 24.400000000000002),
    @(22, 1339, 21),
    @(23, 1345.5, 19),
    @(24, 1350, 16.100000000000001),
    @(25, 1350.5, 14.9),
    @(26, 1348, 11.7),
    @(27, 1357.5, 11.299999999999999),
    @(28, 1364, 7.5),
    @(29, 1364.5, 10.899999999999999),
    @(30, 1358, 12.799999999999999),
    @(31, 1359.5, 15.6),
    @(32, 1370, 21.1),
    @(33, 1367.5, 25.6),
    @(34, 1370, 24.7),
    @(35, 1381.5, 22.400000000000002),
    @(36, 1373, 17.5),
    @(37, 1376.5, 15.1),
    @(38, 1380, 10.899999999999999),
    @(39, 1390.5, 10.4),
    @(40, 1391, 9.5),
    @(41, 1392.5, 10.1),
    @(42, 1387, 18.899999999999999)
)

foreach ($row in $longData) {
    $r = $row[0]
    $long.Cells.Item($r, 1).Value = $row[1]
    $long.Cells.Item($r, 2).Value = $row[2]
}

# ---------------------------------------------------------------
# 3. Selection / active-tab bookkeeping.
#    "medium" ends up with A1:B42 selected (not the active tab).
#    "short" ends up as the active tab with C1:D42 selected - so it
#    must be the LAST sheet activated/selected in this script.
# ---------------------------------------------------------------
$medium.Range("A1:B42").Select()

$short = $wb.Worksheets.Item("short")
$short.Activate()
$short.Range("C1:D42").Select()
